# Generate Report for Handoff
# Adds two newly-handed-off localization files
#   209d7cd6-a6d5-4f88-9714-1b34f51b4240.md
#   83821f61-d375-41db-a637-92d8d1fa8dc4.md
# to the "Overview", "zh-cn" and "de-de" worksheets, pushing the
# previously-last row ("daa0756a-ad85-4048-9b61-3f417e79fbfe.md") down.

$wb = $excel.ActiveWorkbook

$commitSha   = "ba4913602b2baf78110df260af5e2f70b9e19b79"
$repoBase    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# -----------------------------------------------------------------
# Sheet "Overview"
# -----------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

# Row 5 now describes 209d7cd6 (previously described daa0756a)
$wsOv.Range("A5").Value = "209d7cd6-a6d5-4f88-9714-1b34f51b4240.md"
$wsOv.Hyperlinks.Add($wsOv.Range("B5"), ($repoBase + "209d7cd6-a6d5-4f88-9714-1b34f51b4240.md"), "", "", "e2e\209d7cd6-a6d5-4f88-9714-1b34f51b4240.md") | Out-Null
$wsOv.Range("C5").Value = ".md"
$wsOv.Range("D5").Value = ""
$wsOv.Range("E5").Value = "Ready for handoff"
$wsOv.Range("F5").Value = "Ready for handoff"
$wsOv.Range("G5").Value = "2016-08-24 08:42:39"
$wsOv.Range("G5").NumberFormat = $dateFmt

# Row 6 (new) describes 83821f61
$wsOv.Range("A6").Value = "83821f61-d375-41db-a637-92d8d1fa8dc4.md"
$wsOv.Hyperlinks.Add($wsOv.Range("B6"), ($repoBase + "83821f61-d375-41db-a637-92d8d1fa8dc4.md"), "", "", "e2e\83821f61-d375-41db-a637-92d8d1fa8dc4.md") | Out-Null
$wsOv.Range("C6").Value = ".md"
$wsOv.Range("D6").Value = ""
$wsOv.Range("E6").Value = "Ready for handoff"
$wsOv.Range("F6").Value = "Ready for handoff"
$wsOv.Range("G6").Value = "2016-08-24 08:42:39"
$wsOv.Range("G6").NumberFormat = $dateFmt

# Row 7 (new) re-adds daa0756a, which used to be row 5
$wsOv.Range("A7").Value = "daa0756a-ad85-4048-9b61-3f417e79fbfe.md"
$wsOv.Hyperlinks.Add($wsOv.Range("B7"), ($repoBase + "daa0756a-ad85-4048-9b61-3f417e79fbfe.md"), "", "", "e2e\daa0756a-ad85-4048-9b61-3f417e79fbfe.md") | Out-Null
$wsOv.Range("C7").Value = ".md"
$wsOv.Range("D7").Value = ""
$wsOv.Range("E7").Value = "Ready for handoff"
$wsOv.Range("F7").Value = "Ready for handoff"
$wsOv.Range("G7").Value = "2016-08-24 08:40:40"
$wsOv.Range("G7").NumberFormat = $dateFmt

$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G7"))

# -----------------------------------------------------------------
# Sheet "zh-cn"
# -----------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 5 now describes 209d7cd6 (previously described daa0756a)
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), ($repoBase + "209d7cd6-a6d5-4f88-9714-1b34f51b4240.md"), "", "", "209d7cd6-a6d5-4f88-9714-1b34f51b4240.md") | Out-Null
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "False"
$wsZh.Range("G5").Value = "209d7cd6-a6d5-4f88-9714-1b34f51b4240.8f066428eeb3ada2365a68586793ca47a401e97c.zh-cn.xlf"
$wsZh.Range("H5").Value = "2016-08-24 08:42:33"
$wsZh.Range("H5").NumberFormat = $dateFmt
$wsZh.Range("I5").Value = ""
$wsZh.Range("J5").Value = ""
$wsZh.Range("K5").Value = "0001-01-01 00:00:00"
$wsZh.Range("K5").NumberFormat = $dateFmt
$wsZh.Range("L5").Value = ""
$wsZh.Range("M5").Value = "True"
$wsZh.Range("N5").Value = ""
$wsZh.Range("O5").Value = "False"
$wsZh.Range("P5").Value = ""

# Row 6 (new) describes 83821f61
$wsZh.Hyperlinks.Add($wsZh.Range("A6"), ($repoBase + "83821f61-d375-41db-a637-92d8d1fa8dc4.md"), "", "", "83821f61-d375-41db-a637-92d8d1fa8dc4.md") | Out-Null
$wsZh.Range("B6").Value = ".md"
$wsZh.Range("C6").Value = "Ready for handoff"
$wsZh.Range("D6").Value = "e2e"
$wsZh.Range("E6").Value = "ht"
$wsZh.Range("F6").Value = "False"
$wsZh.Range("G6").Value = "83821f61-d375-41db-a637-92d8d1fa8dc4.ee6061f7a2b17e7e69cb2089b0974f08f9352f95.zh-cn.xlf"
$wsZh.Range("H6").Value = "2016-08-24 08:42:33"
$wsZh.Range("H6").NumberFormat = $dateFmt
$wsZh.Range("I6").Value = ""
$wsZh.Range("J6").Value = ""
$wsZh.Range("K6").Value = "0001-01-01 00:00:00"
$wsZh.Range("K6").NumberFormat = $dateFmt
$wsZh.Range("L6").Value = ""
$wsZh.Range("M6").Value = "True"
$wsZh.Range("N6").Value = ""
$wsZh.Range("O6").Value = "False"
$wsZh.Range("P6").Value = ""

# Row 7 (new) re-adds daa0756a, which used to be row 5
$wsZh.Hyperlinks.Add($wsZh.Range("A7"), ($repoBase + "daa0756a-ad85-4048-9b61-3f417e79fbfe.md"), "", "", "daa0756a-ad85-4048-9b61-3f417e79fbfe.md") | Out-Null
$wsZh.Range("B7").Value = ".md"
$wsZh.Range("C7").Value = "Ready for handoff"
$wsZh.Range("D7").Value = "e2e"
$wsZh.Range("E7").Value = "ht"
$wsZh.Range("F7").Value = "False"
$wsZh.Range("G7").Value = "daa0756a-ad85-4048-9b61-3f417e79fbfe.9cc8b8c91d706790aaca959546fb6d1abcb88d8f.zh-cn.xlf"
$wsZh.Range("H7").Value = "2016-08-24 08:40:36"
$wsZh.Range("H7").NumberFormat = $dateFmt
$wsZh.Range("I7").Value = ""
$wsZh.Range("J7").Value = ""
$wsZh.Range("K7").Value = "0001-01-01 00:00:00"
$wsZh.Range("K7").NumberFormat = $dateFmt
$wsZh.Range("L7").Value = ""
$wsZh.Range("M7").Value = "True"
$wsZh.Range("N7").Value = ""
$wsZh.Range("O7").Value = "False"
$wsZh.Range("P7").Value = ""

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P7"))

# -----------------------------------------------------------------
# Sheet "de-de"
# -----------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 5 now describes 209d7cd6 (previously described daa0756a)
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), ($repoBase + "209d7cd6-a6d5-4f88-9714-1b34f51b4240.md"), "", "", "209d7cd6-a6d5-4f88-9714-1b34f51b4240.md") | Out-Null
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "False"
$wsDe.Range("G5").Value = "209d7cd6-a6d5-4f88-9714-1b34f51b4240.8f066428eeb3ada2365a68586793ca47a401e97c.de-de.xlf"
$wsDe.Range("H5").Value = "2016-08-24 08:42:39"
$wsDe.Range("H5").NumberFormat = $dateFmt
$wsDe.Range("I5").Value = ""
$wsDe.Range("J5").Value = ""
$wsDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDe.Range("K5").NumberFormat = $dateFmt
$wsDe.Range("L5").Value = ""
$wsDe.Range("M5").Value = "True"
$wsDe.Range("N5").Value = ""
$wsDe.Range("O5").Value = "False"
$wsDe.Range("P5").Value = ""

# Row 6 (new) describes 83821f61
$wsDe.Hyperlinks.Add($wsDe.Range("A6"), ($repoBase + "83821f61-d375-41db-a637-92d8d1fa8dc4.md"), "", "", "83821f61-d375-41db-a637-92d8d1fa8dc4.md") | Out-Null
$wsDe.Range("B6").Value = ".md"
$wsDe.Range("C6").Value = "Ready for handoff"
$wsDe.Range("D6").Value = "e2e"
$wsDe.Range("E6").Value = "ht"
$wsDe.Range("F6").Value = "False"
$wsDe.Range("G6").Value = "83821f61-d375-41db-a637-92d8d1fa8dc4.ee6061f7a2b17e7e69cb2089b0974f08f9352f95.de-de.xlf"
$wsDe.Range("H6").Value = "2016-08-24 08:42:39"
$wsDe.Range("H6").NumberFormat = $dateFmt
$wsDe.Range("I6").Value = ""
$wsDe.Range("J6").Value = ""
$wsDe.Range("K6").Value = "0001-01-01 00:00:00"
$wsDe.Range("K6").NumberFormat = $dateFmt
$wsDe.Range("L6").Value = ""
$wsDe.Range("M6").Value = "True"
$wsDe.Range("N6").Value = ""
$wsDe.Range("O6").Value = "False"
$wsDe.Range("P6").Value = ""

# Row 7 (new) re-adds daa0756a, which used to be row 5
$wsDe.Hyperlinks.Add($wsDe.Range("A7"), ($repoBase + "daa0756a-ad85-4048-9b61-3f417e79fbfe.md"), "", "", "daa0756a-ad85-4048-9b61-3f417e79fbfe.md") | Out-Null
$wsDe.Range("B7").Value = ".md"
$wsDe.Range("C7").Value = "Ready for handoff"
$wsDe.Range("D7").Value = "e2e"
$wsDe.Range("E7").Value = "ht"
$wsDe.Range("F7").Value = "False"
$wsDe.Range("G7").Value = "daa0756a-ad85-4048-9b61-3f417e79fbfe.9cc8b8c91d706790aaca959546fb6d1abcb88d8f.de-de.xlf"
$wsDe.Range("H7").Value = "2016-08-24 08:40:40"
$wsDe.Range("H7").NumberFormat = $dateFmt
$wsDe.Range("I7").Value = ""
$wsDe.Range("J7").Value = ""
$wsDe.Range("K7").Value = "0001-01-01 00:00:00"
$wsDe.Range("K7").NumberFormat = $dateFmt
$wsDe.Range("L7").Value = ""
$wsDe.Range("M7").Value = "True"
$wsDe.Range("N7").Value = ""
$wsDe.Range("O7").Value = "False"
$wsDe.Range("P7").Value = ""

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P7"))
